$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.02148943381431
$ws.Range("D2").Value = 6.101878592698917
$ws.Range("E2").Value = 13.33796903530186
$ws.Range("F2").Value = 29.46868237896947
$ws.Range("G2").Value = 3.668620872007347
$ws.Range("K2").Value = 15.44481571709794
$ws.Range("L2").Value = 9.17649479592067
$ws.Range("M2").Value = 17.89159082318793
$ws.Range("O2").Value = 26.42479909034651

$ws.Range("C3").Value = 12.99936632806163
$ws.Range("D3").Value = 6.082845908519899
$ws.Range("E3").Value = 13.36868267306908
$ws.Range("F3").Value = 29.57150027655167
$ws.Range("G3").Value = 3.671061570125231
$ws.Range("K3").Value = 14.93775709279146
$ws.Range("L3").Value = 9.203705098353741
$ws.Range("M3").Value = 17.68531480047369
$ws.Range("O3").Value = 26.5564774297373

$ws.Range("C4").Value = 12.98891863063415
$ws.Range("D4").Value = 6.071262935115393
$ws.Range("E4").Value = 13.39015640187051
$ws.Range("F4").Value = 29.64421741816964
$ws.Range("G4").Value = 3.672638510521594
$ws.Range("K4").Value = 14.61845896621865
$ws.Range("L4").Value = 9.221454542092367
$ws.Range("M4").Value = 17.55975735255869
$ws.Range("O4").Value = 26.6446756363319

$ws.Range("C5").Value = 12.98545272501369
$ws.Range("D5").Value = 6.066570462588697
$ws.Range("E5").Value = 13.39956377547744
$ws.Range("F5").Value = 29.67624845918291
$ws.Range("G5").Value = 3.673300891530008
$ws.Range("K5").Value = 14.48652319269797
$ws.Range("L5").Value = 9.228950097597512
$ws.Range("M5").Value = 17.50891720582386
$ws.Range("O5").Value = 26.68245763504008

$ws.Range("C6").Value = 12.98492510091399
$ws.Range("D6").Value = 6.065793003471857
$ws.Range("E6").Value = 13.40116549154802
$ws.Range("F6").Value = 29.68171166107327
$ws.Range("G6").Value = 3.673412075068542
$ws.Range("K6").Value = 14.46451114366066
$ws.Range("L6").Value = 9.230210597082745
$ws.Range("M6").Value = 17.50049634240744
$ws.Range("O6").Value = 26.68884226478974

$ws.Range("C7").Value = 12.98886867949142
$ws.Range("D7").Value = 6.071199536436097
$ws.Range("E7").Value = 13.39028061564526
$ws.Range("F7").Value = 29.64463970660647
$ws.Range("G7").Value = 3.672647363509661
$ws.Range("K7").Value = 14.61668674659325
$ws.Range("L7").Value = 9.221554566284885
$ws.Range("M7").Value = 17.55907032119942
$ws.Range("O7").Value = 26.64517773472397

$ws.Range("C8").Value = 13.01321255211206
$ws.Range("D8").Value = 6.095295757568204
$ws.Range("E8").Value = 13.34801567935016
$ws.Range("F8").Value = 29.50213703537013
$ws.Range("G8").Value = 3.669446204243509
$ws.Range("K8").Value = 15.27173922412856
$ws.Range("L8").Value = 9.185660863366651
$ws.Range("M8").Value = 17.82026915289108
$ws.Range("O8").Value = 26.46867268095908

$ws.Range("C9").Value = 13.08566892957344
$ws.Range("D9").Value = 6.143290898204762
$ws.Range("E9").Value = 13.28592908880002
$ws.Range("F9").Value = 29.29928062963338
$ws.Range("G9").Value = 3.663787364642809
$ws.Range("K9").Value = 16.48589509210368
$ws.Range("L9").Value = 9.123522872450087
$ws.Range("M9").Value = 18.33898446761839
$ws.Range("O9").Value = 26.18115470061339

$ws.Range("C10").Value = 13.15370700262316
$ws.Range("D10").Value = 6.178904121072875
$ws.Range("E10").Value = 13.25304492477678
$ws.Range("F10").Value = 29.19760947095903
$ws.Range("G10").Value = 3.660002761890659
$ws.Range("K10").Value = 17.32645172915916
$ws.Range("L10").Value = 9.082871393195401
$ws.Range("M10").Value = 18.72110813439992
$ws.Range("O10").Value = 26.00606814467543

$ws.Range("C11").Value = 13.1878030325351
$ws.Range("D11").Value = 6.195159848530275
$ws.Range("E11").Value = 13.24085893170866
$ws.Range("F11").Value = 29.16177089849388
$ws.Range("G11").Value = 3.658361138520943
$ws.Range("K11").Value = 17.69615379672766
$ws.Range("L11").Value = 9.065458166705543
$ws.Range("M11").Value = 18.89455086178097
$ws.Range("O11").Value = 25.93435726636874

$ws.Range("C12").Value = 13.20115948494104
$ws.Range("D12").Value = 6.201321474057477
$ws.Range("E12").Value = 13.23664377090328
$ws.Range("F12").Value = 29.14970625691047
$ws.Range("G12").Value = 3.657750936290341
$ws.Range("K12").Value = 17.83422099345435
$ws.Range("L12").Value = 9.059019017421715
$ws.Range("M12").Value = 18.9601227024849
$ws.Range("O12").Value = 25.90835097751964

$ws.Range("C13").Value = 13.19826326364902
$ws.Range("D13").Value = 6.199994225608778
$ws.Range("E13").Value = 13.23753380856822
$ws.Range("F13").Value = 29.15223744494509
$ws.Range("G13").Value = 3.657881846204158
$ws.Range("K13").Value = 17.80457322916952
$ws.Range("L13").Value = 9.060398921331423
$ws.Range("M13").Value = 18.9460062065894
$ws.Range("O13").Value = 25.91390066948408

$ws.Range("C14").Value = 13.18889300076974
$ws.Range("D14").Value = 6.195666661671552
$ws.Range("E14").Value = 13.24050413993236
$ws.Range("F14").Value = 29.16074808715392
$ws.Range("G14").Value = 3.658310707788755
$ws.Range("K14").Value = 17.70755194517431
$ws.Range("L14").Value = 9.064925312310264
$ws.Range("M14").Value = 18.8999479349445
$ws.Range("O14").Value = 25.93219464021896

$ws.Range("C15").Value = 13.18321117341844
$ws.Range("D15").Value = 6.19301661740618
$ws.Range("E15").Value = 13.24237558532602
$ws.Range("F15").Value = 29.16615757446884
$ws.Range("G15").Value = 3.658574886566111
$ws.Range("K15").Value = 17.64786912928149
$ws.Range("L15").Value = 9.067718014926724
$ws.Range("M15").Value = 18.8717204503428
$ws.Range("O15").Value = 25.94355009651089

$ws.Range("C16").Value = 13.15154143148547
$ws.Range("D16").Value = 6.177842710034679
$ws.Range("E16").Value = 13.25389716817354
$ws.Range("F16").Value = 29.20016194874171
$ws.Range("G16").Value = 3.660111650679378
$ws.Range("K16").Value = 17.30202625112026
$ws.Range("L16").Value = 9.084031077773286
$ws.Range("M16").Value = 18.70976102717882
$ws.Range("O16").Value = 26.01091499356663

$ws.Range("C17").Value = 13.13291368007086
$ws.Range("D17").Value = 6.168546709975518
$ws.Range("E17").Value = 13.26167605028697
$ws.Range("F17").Value = 29.22369630244686
$ws.Range("G17").Value = 3.661074855259831
$ws.Range("K17").Value = 17.08653583624413
$ws.Range("L17").Value = 9.094314808811284
$ws.Range("M17").Value = 18.61026802528296
$ws.Range("O17").Value = 26.05427953594017

$ws.Range("C18").Value = 13.12249600496157
$ws.Range("D18").Value = 6.163205202095949
$ws.Range("E18").Value = 13.26641125196905
$ws.Range("F18").Value = 29.23821217497498
$ws.Range("G18").Value = 3.661636399708305
$ws.Range("K18").Value = 16.96140440894613
$ws.Range("L18").Value = 9.100331351189189
$ws.Range("M18").Value = 18.55300885048969
$ws.Range("O18").Value = 26.07996826999555

$ws.Range("C19").Value = 13.11901989575677
$ws.Range("D19").Value = 6.161397630900073
$ws.Range("E19").Value = 13.26805931476657
$ws.Range("F19").Value = 29.2432949350406
$ws.Range("G19").Value = 3.661827824947582
$ws.Range("K19").Value = 16.91883671125193
$ws.Range("L19").Value = 9.102385911730062
$ws.Range("M19").Value = 18.53361772177688
$ws.Range("O19").Value = 26.08879402537965

$ws.Range("C20").Value = 13.13486599822844
$ws.Range("D20").Value = 6.169535750300424
$ws.Range("E20").Value = 13.26082095810851
$ws.Range("F20").Value = 29.2210895811581
$ws.Range("G20").Value = 3.660971541120238
$ws.Range("K20").Value = 17.10959884218116
$ws.Range("L20").Value = 9.093209574119305
$ws.Range("M20").Value = 18.62086306154313
$ws.Range("O20").Value = 26.04958598004975

$ws.Range("C21").Value = 13.19163325724553
$ws.Range("D21").Value = 6.196937626950807
$ws.Range("E21").Value = 13.23962083730166
$ws.Range("F21").Value = 29.15820734491914
$ws.Range("G21").Value = 3.658184430615337
$ws.Range("K21").Value = 17.73610265752917
$ws.Range("L21").Value = 9.063591601370199
$ws.Range("M21").Value = 18.91347967856688
$ws.Range("O21").Value = 25.92679001084565

$ws.Range("C22").Value = 13.23132407503076
$ws.Range("D22").Value = 6.214879996648291
$ws.Range("E22").Value = 13.22809360816133
$ws.Range("F22").Value = 29.12589576788882
$ws.Range("G22").Value = 3.656429574633746
$ws.Range("K22").Value = 18.13426238443354
$ws.Range("L22").Value = 9.045137041754003
$ws.Range("M22").Value = 19.10407393278955
$ws.Range("O22").Value = 25.8532369183984

$ws.Range("C23").Value = 13.20990590346864
$ws.Range("D23").Value = 6.205301396391418
$ws.Range("E23").Value = 13.23403269301064
$ws.Range("F23").Value = 29.14233431717634
$ws.Range("G23").Value = 3.65736009283287
$ws.Range("K23").Value = 17.92282339075573
$ws.Range("L23").Value = 9.054904120995436
$ws.Range("M23").Value = 19.00242615204037
$ws.Range("O23").Value = 25.89187788791441

$ws.Range("C24").Value = 13.13398244614458
$ws.Range("D24").Value = 6.169088596004369
$ws.Range("E24").Value = 13.26120672604616
$ws.Range("F24").Value = 29.22226501023387
$ws.Range("G24").Value = 3.66101822518697
$ws.Range("K24").Value = 17.09917592084382
$ws.Range("L24").Value = 9.093708925837332
$ws.Range("M24").Value = 18.61607322648069
$ws.Range("O24").Value = 26.05170557641761

$ws.Range("C25").Value = 13.06344580234004
$ws.Range("D25").Value = 6.130237839166509
$ws.Range("E25").Value = 13.30049296090898
$ws.Range("F25").Value = 29.34588637952891
$ws.Range("G25").Value = 3.665252436640443
$ws.Range("K25").Value = 16.16591452619996
$ws.Range("L25").Value = 9.139452585919415
$ws.Range("M25").Value = 18.19826655290242
$ws.Range("O25").Value = 26.00606814467543
